$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "The same as experiment in row 11 except here we run it on DLT2 machine for time benchmarking purposes. Currently DLT2 is not as busy as DLT1. This helps because some operations are done on CPU, i.e., DLT1 too busy ==> less CPU cores for my code!"
